# Minor refactoring: split the previously single "General Information" data
# row away from the header by inserting extra blank rows above it, and
# renumber the replicate index / identifier suffix for the CONTROL (water)
# rows on the "Exposure conditions" sheet so they start at 0 instead of 1.

$wb = $excel.ActiveWorkbook

# --- Sheet "General Information": push the single data row from row 2
#     down to row 9 by inserting 7 blank rows above it. ---
$ws1 = $wb.Worksheets.Item("General Information")
$ws1.Rows("2:8").Insert()

# --- Sheet "Exposure conditions": for every CONTROL (water) block
#     (rows 6-9, 14-17, 22-25) decrement the "replicate" number (col I)
#     and the trailing digit of the "PrecisionTox short identifier"
#     (col N) by one. ---
$ws2 = $wb.Worksheets.Item("Exposure conditions")
$controlRows = @(6, 7, 8, 9, 14, 15, 16, 17, 22, 23, 24, 25)

foreach ($r in $controlRows) {
    $cellI = $ws2.Cells.Item($r, 9)
    $cellI.Value = $cellI.Value() - 1

    $cellN = $ws2.Cells.Item($r, 14)
    $identifier = $cellN.Value()
    $prefix = $identifier.Substring(0, $identifier.Length - 1)
    $suffix = [int]$identifier.Substring($identifier.Length - 1) - 1
    $cellN.Value = $prefix + $suffix
}
